$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Equipment")

# Data layout (columns): A=Tag, B=Description, C=PumpType, D=PumpDriverType,
# E=DesignTemp, F=DesignPressure, G=Capacity, H=SpecificGravity, I=DifferentialPressure

$data = @(
    @("Equip-001", "DESC-8", "PT-4", "PDT-5", 4, 4, 6, 5, 8),
    @("Equip-002", "DESC-5", "PT-2", "PDT-4", 3, 8, 4, 3, 7),
    @("Equip-003", "DESC-6", "PT-6", "PDT-5", 6, 7, 4, 6, 3),
    @("Equip-004", "DESC-7", "PT-5", "PDT-4", 4, 4, 2, 8, 9)
)

$row = 2
foreach ($rowValues in $data) {
    $col = 1
    foreach ($val in $rowValues) {
        $ws.Cells.Item($row, $col).Value = $val
        $col++
    }
    $row++
}
